# Auto-generated Excel COM-interop script applying the Marilith_Profits.xlsx data update
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific leve rows
# across all eight sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H58").Value = 1361.8889
$ws.Range("I58").Value = 783.3333
$ws.Range("K58").Value = 2349.9999
$ws.Range("M58").Value = -2199.9999

$ws.Range("H76").Value = 4240.6
$ws.Range("J76").Value = 3750
$ws.Range("L76").Value = 3750
$ws.Range("N76").Value = -4380

$ws.Range("H79").Value = 4240.6
$ws.Range("J79").Value = 3750
$ws.Range("L79").Value = 3750
$ws.Range("N79").Value = -5934

$ws.Range("H88").Value = 886.0833
$ws.Range("J88").Value = 1549.3334
$ws.Range("L88").Value = 1549.3334
$ws.Range("N88").Value = -2361.3334

$ws.Range("H91").Value = 886.0833
$ws.Range("J91").Value = 1549.3334
$ws.Range("L91").Value = 1549.3334
$ws.Range("N91").Value = -4357.3334

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H137").Value = 2263.875
$ws.Range("I137").Value = 2032.5714
$ws.Range("K137").Value = 6097.7142
$ws.Range("M137").Value = -3547.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 36.333332
$ws.Range("I5").Value = 30
$ws.Range("J5").Value = 41.4
$ws.Range("K5").Value = 30
$ws.Range("L5").Value = 41.4
$ws.Range("M5").Value = 82
$ws.Range("N5").Value = -265.4

$ws.Range("H32").Value = 12960
$ws.Range("I32").Value = 11540
$ws.Range("K32").Value = 11540
$ws.Range("M32").Value = -11253

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H61").Value = 1800
$ws.Range("I61").Value = 1800
$ws.Range("K61").Value = 1800
$ws.Range("M61").Value = -1588

$ws.Range("H88").Value = 1982.174
$ws.Range("J88").Value = 2602.375
$ws.Range("L88").Value = 2602.375
$ws.Range("N88").Value = -3414.375

$ws.Range("H91").Value = 1982.174
$ws.Range("J91").Value = 2602.375
$ws.Range("L91").Value = 2602.375
$ws.Range("N91").Value = -5410.375

$ws.Range("H136").Value = 1800
$ws.Range("I136").Value = 1800
$ws.Range("K136").Value = 5400
$ws.Range("M136").Value = -2850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 36.333332
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 41.4
$ws.Range("K4").Value = 30
$ws.Range("L4").Value = 41.4
$ws.Range("M4").Value = 85
$ws.Range("N4").Value = -271.4

$ws.Range("H99").Value = 3087.375
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996

$ws.Range("H105").Value = 1914.1428
$ws.Range("J105").Value = 1883.3334
$ws.Range("L105").Value = 1883.3334
$ws.Range("N105").Value = -5377.3334

$ws.Range("H107").Value = 1273.5
$ws.Range("I107").Value = 1347
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1347
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 573
$ws.Range("N107").Value = -5040

$ws.Range("H110").Value = 64999.5
$ws.Range("J110").Value = 64999.5
$ws.Range("L110").Value = 64999.5
$ws.Range("N110").Value = -73179.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3132.8
$ws.Range("I31").Value = 2059.6
$ws.Range("K31").Value = 2059.6
$ws.Range("M31").Value = -1764.6

$ws.Range("H34").Value = 3132.8
$ws.Range("I34").Value = 2059.6
$ws.Range("K34").Value = 2059.6
$ws.Range("M34").Value = -1857.6

$ws.Range("H58").Value = 2080.6956
$ws.Range("I58").Value = 2066.7896
$ws.Range("K58").Value = 2066.7896
$ws.Range("M58").Value = -1863.7896

$ws.Range("H136").Value = 2080.6956
$ws.Range("I136").Value = 2066.7896
$ws.Range("K136").Value = 6200.3688
$ws.Range("M136").Value = -3650.3688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 177.5
$ws.Range("I12").Value = 148
$ws.Range("J12").Value = 226.66667
$ws.Range("K12").Value = 444
$ws.Range("L12").Value = 680.00001
$ws.Range("M12").Value = -271
$ws.Range("N12").Value = -1026.00001

$ws.Range("H18").Value = 1226
$ws.Range("I18").Value = 1226
$ws.Range("K18").Value = 3678
$ws.Range("M18").Value = -3509

$ws.Range("H101").Value = 5000
$ws.Range("J101").Value = 5000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -19868

$ws.Range("H116").Value = 3498
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 9000
$ws.Range("M116").Value = -5558

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 66000
$ws.Range("I111").Value = 67000
$ws.Range("J111").Value = 65000
$ws.Range("K111").Value = 67000
$ws.Range("L111").Value = 65000
$ws.Range("M111").Value = -63933
$ws.Range("N111").Value = -71134

$ws.Range("H122").Value = 5211093.5
$ws.Range("I122").Value = 5684511
$ws.Range("K122").Value = 17053533
$ws.Range("M122").Value = -17051083

$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -60119

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1379.9
$ws.Range("I22").Value = 1450
$ws.Range("J22").Value = 749
$ws.Range("K22").Value = 1450
$ws.Range("L22").Value = 749
$ws.Range("M22").Value = -1155
$ws.Range("N22").Value = -1339

$ws.Range("H27").Value = 1379.9
$ws.Range("I27").Value = 1450
$ws.Range("J27").Value = 749
$ws.Range("K27").Value = 1450
$ws.Range("L27").Value = 749
$ws.Range("M27").Value = -1343
$ws.Range("N27").Value = -963

$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1451.8
$ws.Range("J107").Value = 1349.5
$ws.Range("L107").Value = 4048.5
$ws.Range("N107").Value = -7888.5

